$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1477832512315271
$ws.Range("C2").Value = 0.5960591133004927
$ws.Range("J2").Value = 0.009852216748768473
$ws.Range("P2").Value = 0.1133004926108374
$ws.Range("S2").Value = 0.1330049261083744
# Row 3
$ws.Range("B3").Value = 0.04580152671755725
$ws.Range("C3").Value = 0.06870229007633588
$ws.Range("J3").Value = 0.01526717557251908
$ws.Range("P3").Value = 0.7404580152671756
$ws.Range("S3").Value = 0.1297709923664122
# Row 4
$ws.Range("J4").Value = 0.06666666666666667
$ws.Range("P4").Value = 0.7
$ws.Range("S4").Value = 0.2333333333333333
# Row 6
$ws.Range("B6").Value = 0.04123711340206185
$ws.Range("E6").Value = 0.005154639175257732
$ws.Range("F6").Value = 0.03608247422680412
$ws.Range("J6").Value = 0.2680412371134021
$ws.Range("O6").Value = 0.03092783505154639
$ws.Range("Q6").Value = 0.1649484536082474
$ws.Range("R6").Value = 0.1030927835051546
$ws.Range("S6").Value = 0.3505154639175257
# Row 7
$ws.Range("B7").Value = 0.1012658227848101
$ws.Range("D7").Value = 0.006329113924050633
$ws.Range("E7").Value = 0.01265822784810127
$ws.Range("F7").Value = 0.03164556962025317
$ws.Range("J7").Value = 0.2278481012658228
$ws.Range("O7").Value = 0.05063291139240506
$ws.Range("Q7").Value = 0.120253164556962
$ws.Range("R7").Value = 0.1075949367088608
$ws.Range("S7").Value = 0.3417721518987342
# Row 8
$ws.Range("B8").Value = 0.08232445520581114
$ws.Range("D8").Value = 0.01210653753026634
$ws.Range("E8").Value = 0.002421307506053269
$ws.Range("F8").Value = 0.05084745762711865
$ws.Range("J8").Value = 0.1331719128329298
$ws.Range("O8").Value = 0.01210653753026634
$ws.Range("Q8").Value = 0.1912832929782082
$ws.Range("R8").Value = 0.1283292978208233
$ws.Range("S8").Value = 0.387409200968523
# Row 9
$ws.Range("B9").Value = 0.1118421052631579
$ws.Range("D9").Value = 0.006578947368421052
$ws.Range("F9").Value = 0.07894736842105263
$ws.Range("J9").Value = 0.1118421052631579
$ws.Range("O9").Value = 0.03289473684210526
$ws.Range("Q9").Value = 0.1973684210526316
$ws.Range("R9").Value = 0.1118421052631579
$ws.Range("S9").Value = 0.3486842105263158
# Row 10
$ws.Range("B10").Value = 0.0753646677471637
$ws.Range("D10").Value = 0.0186385737439222
$ws.Range("E10").Value = 0.002431118314424636
$ws.Range("F10").Value = 0.07212317666126418
$ws.Range("J10").Value = 0.1393841166936791
$ws.Range("O10").Value = 0.01377633711507293
$ws.Range("Q10").Value = 0.2074554294975689
$ws.Range("R10").Value = 0.1353322528363047
$ws.Range("S10").Value = 0.3354943273905997
# Row 11
$ws.Range("G11").Value = 0.1171875
$ws.Range("J11").Value = 0.12109375
$ws.Range("K11").Value = 0.1953125
$ws.Range("L11").Value = 0.52734375
$ws.Range("S11").Value = 0.0390625
# Row 12
$ws.Range("G12").Value = 0.7346938775510204
$ws.Range("J12").Value = 0.163265306122449
$ws.Range("K12").Value = 0.01360544217687075
$ws.Range("L12").Value = 0.03401360544217687
$ws.Range("S12").Value = 0.05442176870748299
# Row 13
$ws.Range("G13").Value = 0.7058823529411765
$ws.Range("J13").Value = 0.2058823529411765
$ws.Range("S13").Value = 0.08823529411764706
# Row 15
$ws.Range("F15").Value = 0.02702702702702703
$ws.Range("H15").Value = 0.1486486486486487
$ws.Range("I15").Value = 0.04954954954954955
$ws.Range("J15").Value = 0.4009009009009009
$ws.Range("K15").Value = 0.07207207207207207
$ws.Range("M15").Value = 0.004504504504504504
$ws.Range("O15").Value = 0.03603603603603604
$ws.Range("S15").Value = 0.2612612612612613
# Row 16
$ws.Range("H16").Value = 0.1714285714285714
$ws.Range("I16").Value = 0.07142857142857142
$ws.Range("J16").Value = 0.4642857142857143
$ws.Range("K16").Value = 0.06428571428571428
$ws.Range("M16").Value = 0.01428571428571429
$ws.Range("O16").Value = 0.05
$ws.Range("S16").Value = 0.1642857142857143
# Row 17
$ws.Range("F17").Value = 0.02386634844868735
$ws.Range("H17").Value = 0.1718377088305489
$ws.Range("I17").Value = 0.08353221957040573
$ws.Range("J17").Value = 0.4534606205250596
$ws.Range("K17").Value = 0.08353221957040573
$ws.Range("M17").Value = 0.02147971360381861
$ws.Range("O17").Value = 0.06682577565632458
$ws.Range("S17").Value = 0.0954653937947494
# Row 18
$ws.Range("F18").Value = 0.01454545454545455
$ws.Range("H18").Value = 0.2109090909090909
$ws.Range("I18").Value = 0.08727272727272728
$ws.Range("J18").Value = 0.4290909090909091
$ws.Range("K18").Value = 0.05454545454545454
$ws.Range("M18").Value = 0.01818181818181818
$ws.Range("N18").Value = 0.003636363636363636
$ws.Range("O18").Value = 0.08
$ws.Range("S18").Value = 0.1018181818181818
# Row 19
$ws.Range("F19").Value = 0.01685393258426966
$ws.Range("H19").Value = 0.2134831460674157
$ws.Range("I19").Value = 0.06647940074906367
$ws.Range("J19").Value = 0.3548689138576779
$ws.Range("K19").Value = 0.1217228464419476
$ws.Range("M19").Value = 0.01779026217228464
$ws.Range("O19").Value = 0.08426966292134831
$ws.Range("S19").Value = 0.1245318352059925
